$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

$ws.Range("J7").Value = "nan"
$ws.Range("J8").Value = "✅"
